# Post_Flight_Datasheet_07_03_V3.xlsx
# Convert the comma-decimal text values (stored as shared strings) in the
# Stationary measurements / Elevator Trim Curve tables into real numeric
# values, so the sheet works with both the "ref" and "true" data sets.
# Also updates the saved view/selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "F28" = 1.8
    "F29" = 2.1
    "J29" = -1.5
    "F30" = 3.4
    "J30" = -3.8
    "F31" = 5.3
    "J31" = -4.5
    "F32" = 8.3
    "H32" = 433.5
    "J32" = -6.9
    "F33" = 10.5
    "J33" = -7.2

    "F59" = 6.2
    "G59" = -0.2
    "H59" = 3.7
    "I59" = -35.5

    "F60" = 5.4
    "G60" = 0.2
    "H60" = 3.7
    "M60" = -5.5

    "F61" = 4.6
    "G61" = 0.5
    "H61" = 3.7
    "M61" = -4.8

    "F62" = 3.7
    "G62" = 0.9
    "H62" = 3.7
    "M62" = -2.8

    "F63" = 3.4
    "G63" = 1.1
    "H63" = 3.7
    "I63" = 53.5
    "M63" = -2.2

    "F75" = 4.7
    "G75" = 0.4
    "H75" = 3.7
    "M75" = -4.8

    "F76" = 4.6
    "G76" = 0.1
    "H76" = 3.7
    "M76" = -4.5
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# Update the stored view/selection state for the sheet
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 64
$win.ScrollColumn = 1
$ws.Range("N28").Select()
